# Update dSF (column F) values on Sheet1 to reflect repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F5").Value = 5
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 2
$ws.Range("F13").Value = 2
$ws.Range("F19").Value = 2
$ws.Range("F25").Value = 3
$ws.Range("F29").Value = -2
$ws.Range("F31").Value = -6
$ws.Range("F35").Value = -6
$ws.Range("F51").Value = -1
$ws.Range("F55").Value = -2
$ws.Range("F57").Value = 0
$ws.Range("F60").Value = -4
$ws.Range("F61").Value = 1
$ws.Range("F62").Value = 1
$ws.Range("F65").Value = 1
$ws.Range("F66").Value = -7
$ws.Range("F69").Value = -1
$ws.Range("F70").Value = 1
$ws.Range("F72").Value = 5
